$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 data updates (Reg iProctor TC user details)
$ws.Range("A2").Value = "AgGbg661"
$ws.Range("B2").Value = 23073124
$ws.Range("C2").Value = "gvxdihi52"
$ws.Range("D2").Value = "Qf#8$9mA"
$ws.Range("F2").Value = "AAKFrBgj"
$ws.Range("G2").Value = "idAs"
